# Add a "number_of_run" parameter column to the scenarios sheet and tidy up
# the active sheet/selection so the "scenarios" sheet is shown with the new
# column selected.

$wb = $excel.ActiveWorkbook
$wsScenarios = $wb.Worksheets.Item("scenarios")
$wsZero = $wb.Worksheets.Item("0")

# Insert a new column before the old column D ("agent_account_min"),
# shifting the rest of the parameter columns one place to the right.
$wsScenarios.Range("D1").EntireColumn.Insert()

# Header + values for the new "number_of_run" parameter column.
$wsScenarios.Range("D1").Value2 = "number_of_run"
$wsScenarios.Range("D2").Value2 = 1
$wsScenarios.Range("D3").Value2 = 2
$wsScenarios.Range("D4").Value2 = 1

# Restore/refresh the selection on the "0" sheet (it was previously the
# active tab) before switching focus away from it.
$wsZero.Range("C102").Select()

# Make "scenarios" the active sheet again and update its selection to the
# newly inserted column.
$wsScenarios.Select()
$wsScenarios.Range("D5").Select()
